# Thinh comimit vs 25.12.2019
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4, 6, 8 are "Buy off sample" / NG rows that previously had
# "T" in the "IS REPAIRED" column (L). Clear that marker.
$ws.Range("L4").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("L8").ClearContents()

# Row 8 also had a stray boolean TRUE in "IS NOT SHOWN" (M) - clear it.
$ws.Range("M8").ClearContents()

# Rows 5, 7, 9 are the corresponding "Repair" / OK rows. They now get
# "F" recorded in both L (IS REPAIRED) and M (IS NOT SHOWN).
$ws.Range("L5").Value = "F"
$ws.Range("M5").Value = "F"

$ws.Range("L7").Value = "F"
$ws.Range("M7").Value = "F"

$ws.Range("L9").Value = "F"
$ws.Range("M9").Value = "F"
